$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the PHPSESSID query-string value embedded in the Product Url column
# for the two remaining data rows.
$ws.Range("K2").Value = "https://www.leguidedesmontres.com/en/products-new/a-lange-amp-sohne/1003/233021?PHPSESSID=856f944be12372379eaf709910c4a801"
$ws.Range("K3").Value = "https://www.leguidedesmontres.com/en/products-new/a-lange-amp-sohne/1001/233025?PHPSESSID=856f944be12372379eaf709910c4a801"

# Remove the trailing rows (4-11) that only contained stale/empty records,
# shrinking the used range down to A1:K3.
$ws.Range("A4:K11").Delete()
